$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.875.53"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.841.39"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'309.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.4701"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.45%  "
$ws.Range("D8").Value = "'0.3659"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "'0.07141"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'0.9216"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'19.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07659"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.768.11"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "'5.287"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'6.397"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "'88.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'0.000008649"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "26.904.76"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'5.008"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "'1.930"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'151.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'18.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "'2.008"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'114.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "'4.878"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").Value = "'0.08821"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").Value = "'3.213"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").Value = "'1.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.31%  "
$ws.Range("D33").Value = "'0.7471"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.482"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'2.751"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").Value = "'1.086"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'0.01941"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "'0.05211"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").Value = "'2.958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").Value = "'0.5203"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Value = "'6.976"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").Value = "'0.1511"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "'8.162"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "'10.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("D45").Value = "'0.4698"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "'1.006"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "'102.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("D48").Value = "'1.598"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("D49").Value = "'66.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("D50").Value = "'0.06038"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'0.8879"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.04%  "
